# Generate Report for Handoff
# The handoff run regenerated the localization status report under a new
# source doc id (88005e44-... -> eb7f4660-...) and refreshed the handoff
# timestamps / xliff file names that depend on it.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldId = "88005e44-57c9-4227-b6c3-cbd86228ab60"
$newId = "eb7f4660-e7f0-479b-970e-ddae80723102"

$newFileName   = "$newId.md"
$newPathName   = "e2e\$newId.md"
$newZhXlf      = "$newId.02d2ff92687eb9593844a112278261a058fbddde.zh-cn.xlf"
$newDeXlf      = "$newId.02d2ff92687eb9593844a112278261a058fbddde.de-de.xlf"

$newHoDate     = "2016-08-27 08:57:05"
$newZhHoDate   = "2016-08-27 08:56:57"
$newDeHoDate   = "2016-08-27 08:57:05"

# --- Overview sheet: File Name / Path And Name / Latest HO Xliff Generate Date ---
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = $newHoDate

# --- zh-cn sheet: Source File Name / Latest Handoff File / Latest Handoff Datetime ---
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhHoDate

# --- de-de sheet: Source File Name / Latest Handoff File / Latest Handoff Datetime ---
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newDeHoDate

# --- Update the hyperlink display text on each sheet to match the new file name.
#     Iterating with foreach (rather than Hyperlinks.Item(1)) updates the
#     existing link's display text in place instead of inserting a duplicate.
foreach ($hl in $wsOverview.Hyperlinks) { $hl.TextToDisplay = $newPathName }
foreach ($hl in $wsZhCn.Hyperlinks)     { $hl.TextToDisplay = $newFileName }
foreach ($hl in $wsDeDe.Hyperlinks)     { $hl.TextToDisplay = $newFileName }
